$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire rows (in descending order so row indices of rows still to
# delete remain valid as rows above get removed).
$ws.Rows.Item(40).Delete()
$ws.Rows.Item(37).Delete()
$ws.Rows.Item(36).Delete()
$ws.Rows.Item(12).Delete()
